$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value of B4 to the new text (adds a new shared string entry)
$ws.Range("B4").Value = "2 timer 30 minutter"

# Update the active selection from B10 to B6, matching the saved view state
$ws.Range("B6").Select()
